$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.235.85'
$ws.Range('D3').Value = '1.783.56'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '338.02'
$ws.Range('E5').Value = '  +0.35%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3830'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3445'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.11'
$ws.Range('E9').Value = '  -1.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.157'
$ws.Range('E10').Value = '  -2.97%  '
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '23.42'
$ws.Range('E12').Value = '  +7.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.001'
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.438'
$ws.Range('E14').Value = '  -0.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.220'
$ws.Range('D16').Value = '1.785.22'
$ws.Range('E16').Value = '  +0.20%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001076'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06667'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '82.64'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9999'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.53'
$ws.Range('E21').Value = '  +0.40%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.465'
$ws.Range('E22').Value = '  -0.99%  '
$ws.Range('D23').Value = '28.227.99'
$ws.Range('E23').Value = '  +3.98%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.12'
$ws.Range('E24').Value = '  -1.58%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.368'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('B26').Value = 'ImmutableX'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.445'
$ws.Range('E26').Value = '  +0.19%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.85'
$ws.Range('E27').Value = '  -1.59%  '
$ws.Range('B28').Value = 'Monero'
$ws.Range('C28').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.95'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.426'
$ws.Range('E29').Value = '  -3.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '137.66'
$ws.Range('E30').Value = '  +2.39%  '
$ws.Range('D31').Value = '1.988.80'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.161'
$ws.Range('E32').Value = '  +2.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.982'
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08906'
$ws.Range('E34').Value = '  +2.61%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.84'
$ws.Range('E35').Value = '  -1.80%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02431'
$ws.Range('E36').Value = '  +3.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.6872'
$ws.Range('E37').Value = '  +0.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.343'
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.06362'
$ws.Range('E39').Value = '  +0.50%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.241'
$ws.Range('E41').Value = '  +0.20%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.498'
$ws.Range('E42').Value = '  -7.60%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.315'
$ws.Range('E43').Value = '  -1.55%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.28'
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9993'
$ws.Range('E45').Value = '  -0.07%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6306'
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.870'
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.29'
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.097'
$ws.Range('E49').Value = '  -3.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07488'
$ws.Range('E50').Value = '  +5.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.218'
$ws.Range('E51').Value = '  +9.34%  '
